# Add a new "ODI Bowling Extra" worksheet (mirrors the existing
# "ODI Batting Extra" sheet) containing per-match bowling extras:
# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL.

$wb = $excel.ActiveWorkbook

# --- Create the new sheet and park it at the end of the tab strip ---
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "ODI Bowling Extra"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Move($null, $lastSheet)

# NOTE: worksheet object references captured before a Move() call go
# stale (they keep tracking the ORIGINAL slot index, not the sheet
# object). Always re-fetch the sheet by name after moving it.
$ws = $wb.Worksheets.Item("ODI Bowling Extra")

# --- Header row ---
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
}

$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL ---
# All three columns are plain text (even the numeric-looking ones),
# matching the rest of the workbook's "_Extra" sheets.
$data = @(
    @("3878", "", ""),
    @("3893", "", ""),
    @("3894", "0", "10.00%"),
    @("3924", "0", "10.00%"),
    @("4169", "", ""),
    @("4170", "0", ""),
    @("4234", "1", ""),
    @("4235", "", ""),
    @("4263", "0", "10.00%"),
    @("4266", "0", "10.00%"),
    @("4270", "1", "10.00%"),
    @("4273", "0", "10.00%"),
    @("4274", "0", "10.00%"),
    @("4275", "0", "10.00%"),
    @("4276", "", ""),
    @("4277", "0", "10.00%"),
    @("4336", "0", ""),
    @("4341", "0", "10.00%"),
    @("4351", "0", "20.00%"),
    @("4354", "0", "")
)

$rowIndex = 2
foreach ($record in $data) {
    $rowRange = $ws.Range("A$rowIndex`:C$rowIndex")
    $rowRange.NumberFormat = "@"

    $ws.Cells.Item($rowIndex, 1).Value = $record[0]
    $ws.Cells.Item($rowIndex, 2).Value = $record[1]
    $ws.Cells.Item($rowIndex, 3).Value = $record[2]

    $rowIndex++
}
